$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.824.89'
$ws.Range("E2").Value = '  +2.61%  '

$ws.Range("D3").Value = '3.552.43'
$ws.Range("E3").Value = '  +1.32%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '610.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.38%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.97'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.46%  '

$ws.Range("E7").Value = '  +1.74%  '

$ws.Range("D8").Value = '3.546.32'
$ws.Range("E8").Value = '  +1.38%  '

$ws.Range("E9").Value = '  -0.03%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.197'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.92'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.54%  '

$ws.Range("E12").Value = '  +0.33%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.72'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.73%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000277'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.88%  '

$ws.Range("D15").Value = '4.124.04'
$ws.Range("E15").Value = '  +1.41%  '

$ws.Range("E16").Value = '  -1.65%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '618.31'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.20%  '

$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '70.826.18'
$ws.Range("E18").Value = '  +2.49%  '

$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '3.549.92'
$ws.Range("E19").Value = '  +0.95%  '

$ws.Range("E20").Value = '  -1.15%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.40'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.52%  '

$ws.Range("E22").Value = '  -0.64%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.48'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -14.57%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.72'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.37%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '96.77'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.91%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.81'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.24%  '

$ws.Range("E27").Value = '  -0.02%  '

$ws.Range("E28").Value = '  -0.93%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.54'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.71%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.06'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.49'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.31%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.08'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.77%  '

$ws.Range("E33").Value = '  -0.78%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.00'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.33%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '573.75'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.61%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.101'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.74%  '

$ws.Range("B37").Value = 'dogwifhat'
$ws.Range("C37").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.61'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.88%  '

$ws.Range("E38").Value = '  +0.40%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '57.67'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.57%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0472'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.46%  '

$ws.Range("E41").Value = '  +0.03%  '

$ws.Range("E42").Value = '  +4.38%  '

$ws.Range("D43").Value = '3.357.13'
$ws.Range("E43").Value = '  -0.11%  '

$ws.Range("E44").Value = '  -2.22%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.80%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '33.01'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.60%  '

$ws.Range("D47").Value = '0.0₃0706'
$ws.Range("E47").Value = '  +1.26%  '

$ws.Range("E48").Value = '  +2.71%  '

$ws.Range("E49").Value = '  -0.06%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.72'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.22%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.68'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.50%  '
